$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1709.07
$ws.Range("I15").Value = 1709.07
$ws.Range("K15").Value = 5127.21
$ws.Range("M15").Value = -4958.21
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H28").Value = 299.46667
$ws.Range("I28").Value = 188.3077
$ws.Range("J28").Value = 1022
$ws.Range("K28").Value = 188.3077
$ws.Range("L28").Value = 1022
$ws.Range("M28").Value = 296.6923
$ws.Range("N28").Value = -1992
$ws.Range("H33").Value = 167.95
$ws.Range("I33").Value = 175.63158
$ws.Range("K33").Value = 175.63158
$ws.Range("M33").Value = 53.36841999999999
$ws.Range("H40").Value = 1318.3572
$ws.Range("I40").Value = 723.3333
$ws.Range("J40").Value = 1764.625
$ws.Range("K40").Value = 723.3333
$ws.Range("L40").Value = 1764.625
$ws.Range("M40").Value = -548.3333
$ws.Range("N40").Value = -2114.625
$ws.Range("H53").Value = 4599.7144
$ws.Range("J53").Value = 6339.4
$ws.Range("L53").Value = 6339.4
$ws.Range("N53").Value = -7613.4
$ws.Range("H62").Value = 2927.1428
$ws.Range("I62").Value = 2655.5557
$ws.Range("J62").Value = 3130.8333
$ws.Range("K62").Value = 2655.5557
$ws.Range("L62").Value = 3130.8333
$ws.Range("M62").Value = -2031.5557
$ws.Range("N62").Value = -4378.8333
$ws.Range("H65").Value = 2927.1428
$ws.Range("I65").Value = 2655.5557
$ws.Range("J65").Value = 3130.8333
$ws.Range("K65").Value = 13277.7785
$ws.Range("L65").Value = 15654.1665
$ws.Range("M65").Value = -10157.7785
$ws.Range("N65").Value = -21894.1665
$ws.Range("H116").Value = 3976.647
$ws.Range("I116").Value = 2666.6667
$ws.Range("K116").Value = 2666.6667
$ws.Range("M116").Value = 775.3332999999998
$ws.Range("H125").Value = 548.0833
$ws.Range("J125").Value = 517.6667
$ws.Range("L125").Value = 4659.0003
$ws.Range("N125").Value = -9579.0003
$ws.Range("H129").Value = 755.87177
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 802.25714
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 2406.77142
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -12406.77142
$ws.Range("H132").Value = 2808.3428
$ws.Range("I132").Value = 2887.6365
$ws.Range("K132").Value = 8662.9095
$ws.Range("M132").Value = -6132.9095
$ws.Range("H137").Value = 1537.037
$ws.Range("I137").Value = 1629.579
$ws.Range("J137").Value = 1317.25
$ws.Range("K137").Value = 4888.737
$ws.Range("L137").Value = 3951.75
$ws.Range("M137").Value = -2338.737
$ws.Range("N137").Value = -9051.75
$ws.Range("H138").Value = 2600.3428
$ws.Range("I138").Value = 1532.2941
$ws.Range("J138").Value = 3609.0557
$ws.Range("K138").Value = 4596.8823
$ws.Range("L138").Value = 10827.1671
$ws.Range("M138").Value = 543.1176999999998
$ws.Range("N138").Value = -21107.1671

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1008.6857
$ws.Range("I122").Value = 1029.5294
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 3088.5882
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = -638.5881999999997
$ws.Range("N122").Value = -5800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1335
$ws.Range("I107").Value = 1180.9
$ws.Range("J107").Value = 1555.1428
$ws.Range("K107").Value = 1180.9
$ws.Range("L107").Value = 1555.1428
$ws.Range("M107").Value = 739.0999999999999
$ws.Range("N107").Value = -5395.1428

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2645.28
$ws.Range("I132").Value = 1528.8889
$ws.Range("J132").Value = 5516
$ws.Range("K132").Value = 4586.6667
$ws.Range("L132").Value = 16548
$ws.Range("M132").Value = -2056.6667
$ws.Range("N132").Value = -21608
$ws.Range("H141").Value = 29990.926
$ws.Range("J141").Value = 29990.926
$ws.Range("L141").Value = 29990.926
$ws.Range("N141").Value = -40350.926

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 35.714287
$ws.Range("I11").Value = 34
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 102
$ws.Range("L11").Value = 120
$ws.Range("M11").Value = 38
$ws.Range("N11").Value = -400
$ws.Range("H87").Value = 21809.523
$ws.Range("J87").Value = 26446.666
$ws.Range("L87").Value = 79339.99800000001
$ws.Range("N87").Value = -81835.99800000001
$ws.Range("H90").Value = 21809.523
$ws.Range("J90").Value = 26446.666
$ws.Range("L90").Value = 238019.994
$ws.Range("N90").Value = -250499.994
$ws.Range("H93").Value = 10000
$ws.Range("J93").Value = 10000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H122").Value = 1234.4783
$ws.Range("I122").Value = 336
$ws.Range("J122").Value = 1423.6316
$ws.Range("K122").Value = 3024
$ws.Range("L122").Value = 12812.6844
$ws.Range("M122").Value = -574
$ws.Range("N122").Value = -17712.6844
$ws.Range("H131").Value = 720.6
$ws.Range("J131").Value = 727.8372000000001
$ws.Range("L131").Value = 2183.5116
$ws.Range("N131").Value = -12263.5116

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4176973.5
$ws.Range("I70").Value = 4401.143
$ws.Range("J70").Value = 7827974.5
$ws.Range("K70").Value = 4401.143
$ws.Range("L70").Value = 7827974.5
$ws.Range("M70").Value = -4131.143
$ws.Range("N70").Value = -7828514.5
$ws.Range("H73").Value = 4176973.5
$ws.Range("I73").Value = 4401.143
$ws.Range("J73").Value = 7827974.5
$ws.Range("K73").Value = 4401.143
$ws.Range("L73").Value = 7827974.5
$ws.Range("M73").Value = -3465.143
$ws.Range("N73").Value = -7829846.5
$ws.Range("H107").Value = 426.44446
$ws.Range("J107").Value = 621
$ws.Range("L107").Value = 621
$ws.Range("N107").Value = -4461

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3194.375
$ws.Range("I61").Value = 1233.3
$ws.Range("K61").Value = 1233.3
$ws.Range("M61").Value = -1031.3
$ws.Range("H113").Value = 3194.375
$ws.Range("I113").Value = 1233.3
$ws.Range("K113").Value = 1233.3
$ws.Range("M113").Value = 936.7
$ws.Range("H122").Value = 702818.6
$ws.Range("I122").Value = 855017.9399999999
$ws.Range("J122").Value = 2702
$ws.Range("K122").Value = 2565053.82
$ws.Range("L122").Value = 8106
$ws.Range("M122").Value = -2562603.82
$ws.Range("N122").Value = -13006
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -34860

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 76934200
$ws.Range("J2").Value = 13529.2
$ws.Range("L2").Value = 13529.2
$ws.Range("N2").Value = -13753.2
$ws.Range("H70").Value = 6015748.5
$ws.Range("J70").Value = 6015748.5
$ws.Range("L70").Value = 6015748.5
$ws.Range("N70").Value = -6016378.5
$ws.Range("H73").Value = 6015748.5
$ws.Range("J73").Value = 6015748.5
$ws.Range("L73").Value = 6015748.5
$ws.Range("N73").Value = -6017932.5
$ws.Range("H107").Value = 45455096
$ws.Range("I107").Value = 83333576
$ws.Range("J107").Value = 918.9
$ws.Range("K107").Value = 250000728
$ws.Range("L107").Value = 2756.7
$ws.Range("M107").Value = -249998808
$ws.Range("N107").Value = -6596.7
$ws.Range("H122").Value = 1704.5883
$ws.Range("I122").Value = 1676.7142
$ws.Range("J122").Value = 1834.6666
$ws.Range("K122").Value = 5030.142599999999
$ws.Range("L122").Value = 5503.9998
$ws.Range("M122").Value = -2580.142599999999
$ws.Range("N122").Value = -10403.9998
$ws.Range("H126").Value = 1368.9714
$ws.Range("I126").Value = 1013.8333
$ws.Range("J126").Value = 3499.8
$ws.Range("K126").Value = 3041.4999
$ws.Range("L126").Value = 10499.4
$ws.Range("M126").Value = -571.4998999999998
$ws.Range("N126").Value = -15439.4
$ws.Range("H132").Value = 1063.2046
$ws.Range("I132").Value = 822
$ws.Range("J132").Value = 1529.5333
$ws.Range("K132").Value = 2466
$ws.Range("L132").Value = 4588.5999
$ws.Range("M132").Value = 64
$ws.Range("N132").Value = -9648.599900000001
